$wb = $excel.ActiveWorkbook

# --- Update rate summary text on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.65 = 14029.3 pesos`n✅ 14029.3 pesos = 3.64 = 949.95 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$cell.Value = $newText

# --- Update numeric rate cells on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 273.75
$ws2.Range("O10").Value = 3840.52
$ws2.Range("N12").Value = 3859
$ws2.Range("O12").Value = 261.3
